# Update the RW forms PDF validation
# Mark SmokeTest = "Yes" for all RW probate-form test rows (rows 26-235)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioMapping")
$ws.Activate()

$ws.Range("D26:D235").Value = "Yes"

# Update the saved window view: scroll so row 233 is the top-left visible
# row, and leave the active selection on D240.
$excel.Goto($ws.Range("A233"), $true)
$ws.Range("D240").Select()
